$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Duplicate the formatting (cell styles) of the last existing row (77)
# onto the two new rows, exactly like copying the row down.
$ws.Range("A77:F77").Copy()
$ws.Range("A78:F78").PasteSpecial(-4122)
$ws.Range("A79:F79").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New data rows for 2025-09-08 (Excel serial date 45908)
# Row 78: 四方坪站
$ws.Cells.Item(78, 1).Value = 45908
$ws.Cells.Item(78, 2).Value = "四方坪站"
$ws.Cells.Item(78, 3).Value = 11341.76
$ws.Cells.Item(78, 4).Value = 9140.48
$ws.Cells.Item(78, 5).Value = 3953.65
$ws.Cells.Item(78, 6).Value = 458

# Row 79: 高岭站
$ws.Cells.Item(79, 1).Value = 45908
$ws.Cells.Item(79, 2).Value = "高岭站"
$ws.Cells.Item(79, 3).Value = 4996.2
$ws.Cells.Item(79, 4).Value = 3781.68
$ws.Cells.Item(79, 5).Value = 1333.17
$ws.Cells.Item(79, 6).Value = 169

# Update the visible scroll position / selection like the original edit
$excel.ActiveWindow.ScrollRow = 67
$ws.Range("H75").Select()
